$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear interest values that are not yet evaluated for this user / rows
$ws.Range("C4:F4").ClearContents()
$ws.Range("C5:F5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()

# Update selection to reflect the cell currently being worked on
$ws.Range("C4:C5").Select()
